$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.45491533333333
$ws.Range("H2").Value = 31.364746
$ws.Range("I2").Value = 0.0134573334963438
$ws.Range("J2").Value = 0.0134573334963438
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.214110666666667
$ws.Range("N2").Value = 21.642332
$ws.Range("O2").Value = 0.4688823795981188
$ws.Range("P2").Value = 0.4688823795981188
$ws.Range("Q2").Value = 75.42291622529689
$ws.Range("R2").Value = 678.806246027672
$ws.Range("S2").Value = 0.006309906552811154
$ws.Range("T2").Value = 0.006309906552811153

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.45491533333333
$ws.Range("H3").Value = 31.364746
$ws.Range("I3").Value = 0.0134573334963438
$ws.Range("J3").Value = 0.0134573334963438
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.110350666666666
$ws.Range("N3").Value = 21.331052
$ws.Range("O3").Value = 0.4621384803214003
$ws.Range("P3").Value = 0.4621384803214003
$ws.Range("Q3").Value = 74.33811421031022
$ws.Range("R3").Value = 669.043027892792
$ws.Range("S3").Value = 0.006219151651178602
$ws.Range("T3").Value = 0.006219151651178601

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.45491533333333
$ws.Range("H4").Value = 31.364746
$ws.Range("I4").Value = 0.0134573334963438
$ws.Range("J4").Value = 0.0134573334963438
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.061296333333333
$ws.Range("N4").Value = 3.183889
$ws.Range("O4").Value = 0.06897914008048092
$ws.Range("P4").Value = 0.06897914008048092
$ws.Range("Q4").Value = 11.09576330857711
$ws.Range("R4").Value = 99.86186977719399
$ws.Range("S4").Value = 0.0009282752923540471
$ws.Range("T4").Value = 0.000928275292354047

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 735.4993083333334
$ws.Range("H5").Value = 2206.497925
$ws.Range("I5").Value = 0.9467182815928301
$ws.Range("J5").Value = 0.9467182815928301
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.214110666666667
$ws.Range("N5").Value = 21.642332
$ws.Range("O5").Value = 0.4688823795981188
$ws.Range("P5").Value = 0.4688823795981188
$ws.Range("Q5").Value = 5305.973405573456
$ws.Range("R5").Value = 47753.7606501611
$ws.Range("S5").Value = 0.4438995206822882
$ws.Range("T5").Value = 0.4438995206822882

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 735.4993083333334
$ws.Range("H6").Value = 2206.497925
$ws.Range("I6").Value = 0.9467182815928301
$ws.Range("J6").Value = 0.9467182815928301
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.110350666666666
$ws.Range("N6").Value = 21.331052
$ws.Range("O6").Value = 0.4621384803214003
$ws.Range("P6").Value = 0.4621384803214003
$ws.Range("Q6").Value = 5229.657997340789
$ws.Range("R6").Value = 47066.9219760671
$ws.Range("S6").Value = 0.437514947947798
$ws.Range("T6").Value = 0.437514947947798

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 735.4993083333334
$ws.Range("H7").Value = 2206.497925
$ws.Range("I7").Value = 0.9467182815928301
$ws.Range("J7").Value = 0.9467182815928301
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.061296333333333
$ws.Range("N7").Value = 3.183889
$ws.Range("O7").Value = 0.06897914008048092
$ws.Range("P7").Value = 0.06897914008048092
$ws.Range("Q7").Value = 780.5827191033694
$ws.Range("R7").Value = 7025.244471930325
$ws.Range("S7").Value = 0.06530381296274401
$ws.Range("T7").Value = 0.06530381296274401

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 30.939307
$ws.Range("H8").Value = 92.81792100000001
$ws.Range("I8").Value = 0.03982438491082609
$ws.Range("J8").Value = 0.03982438491082609
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.214110666666667
$ws.Range("N8").Value = 21.642332
$ws.Range("O8").Value = 0.4688823795981188
$ws.Range("P8").Value = 0.4688823795981188
$ws.Range("Q8").Value = 223.1995846479747
$ws.Range("R8").Value = 2008.796261831772
$ws.Range("S8").Value = 0.01867295236301955
$ws.Range("T8").Value = 0.01867295236301955

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 30.939307
$ws.Range("H9").Value = 92.81792100000001
$ws.Range("I9").Value = 0.03982438491082609
$ws.Range("J9").Value = 0.03982438491082609
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.110350666666666
$ws.Range("N9").Value = 21.331052
$ws.Range("O9").Value = 0.4621384803214003
$ws.Range("P9").Value = 0.4621384803214003
$ws.Range("Q9").Value = 219.9893221536547
$ws.Range("R9").Value = 1979.903899382892
$ws.Range("S9").Value = 0.01840438072242367
$ws.Range("T9").Value = 0.01840438072242367

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 30.939307
$ws.Range("H10").Value = 92.81792100000001
$ws.Range("I10").Value = 0.03982438491082609
$ws.Range("J10").Value = 0.03982438491082609
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.061296333333333
$ws.Range("N10").Value = 3.183889
$ws.Range("O10").Value = 0.06897914008048092
$ws.Range("P10").Value = 0.06897914008048092
$ws.Range("Q10").Value = 32.83577307497433
$ws.Range("R10").Value = 295.521957674769
$ws.Range("S10").Value = 0.002747051825382863
$ws.Range("T10").Value = 0.002747051825382863
